$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.953.16"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.640.46"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'217.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D10").Value = "'20.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.869.09"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.642.90"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D16").Value = "'67.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "26.936.11"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'219.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").Value = "'4.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'9.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'147.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "'15.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "1.267.67"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").Value = "'0.538"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").Value = "'0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'5.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").Value = "1.779.87"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'62.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").Value = "'92.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +18.24%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "'7.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  -0.34%  "
